$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.076.18'
$ws.Range("E2").Value = '  +2.82%  '
$ws.Range("D3").Value = '2.952.83'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.89'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.17'
$ws.Range("E6").Value = '  +2.41%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '2.952.71'
$ws.Range("E8").Value = '  +0.86%  '
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.15'
$ws.Range("E10").Value = '  +3.15%  '
$ws.Range("E11").Value = '  +6.57%  '
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("E13").Value = '  +5.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.77'
$ws.Range("E14").Value = '  -2.05%  '
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("D16").Value = '3.444.57'
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("D17").Value = '63.050.27'
$ws.Range("E17").Value = '  +2.87%  '
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").Value = '2.957.74'
$ws.Range("E19").Value = '  +1.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '442.74'
$ws.Range("E20").Value = '  +2.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.50'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.668'
$ws.Range("E22").Value = '  -1.02%  '
$ws.Range("E23").Value = '  -1.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.00'
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.12'
$ws.Range("E25").Value = '  +2.11%  '
$ws.Range("E26").Value = '  -2.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.77'
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.31'
$ws.Range("E29").Value = '  +6.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.20'
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("E32").Value = '  +15.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.49'
$ws.Range("E33").Value = '  -0.60%  '
$ws.Range("E34").Value = '  -1.06%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.989'
$ws.Range("E36").Value = '  -1.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.09'
$ws.Range("E37").Value = '  +3.94%  '
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.04'
$ws.Range("E39").Value = '  +2.35%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.67'
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.51'
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("E42").Value = '  -4.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.280'
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.80'
$ws.Range("E44").Value = '  -8.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '135.33'
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("D46").Value = '2.693.12'
$ws.Range("E46").Value = '  -0.19%  '
$ws.Range("E47").Value = '  -1.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '360.18'
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.75'
$ws.Range("E51").Value = '  -3.39%  '
